$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# The sheet currently has (rows 144-149):
#   144: otf-cjklee-h | Identity-H | :1:Klee.ttc
#   145: otf-cjklee-v | Identity-V | :1:Klee.ttc
#   146: up-min-l-plane02 | unicode | ZhongHuaSongPlane02-Regular.ttf |
#   147: ut-min-l-plane02 | unicode | ZhongHuaSongPlane02-Regular.ttf | -w 1
#   148: up-min-r-plane02 | unicode | ZhongHuaSongPlane02-Regular.ttf |
#   149: ut-min-r-plane02 | unicode | ZhongHuaSongPlane02-Regular.ttf | -w 1
#
# The target layout re-orders these so the "plane02" block (currently
# 146-149) comes first (144-147), followed by the "cjklee" block (currently
# 144-145) at 148-149, and then 14 brand-new rows of font-substitution data
# are appended at 150-163.
# ---------------------------------------------------------------------------

# Step 1: insert 4 blank rows at the top of the block (144:147). Excel
# inherits the row-144-above formatting for the new blank cells, which keeps
# the existing "s=1" style used by the plane02 rows, and pushes the current
# 144-149 block down to 148-153.
$ws.Rows("144:147").Insert(-4121)

# Step 2: the former rows 146-149 (the plane02 block) are now sitting at
# 150-153; copy their values up into the freshly-inserted 144-147 block.
for ($i = 0; $i -lt 4; $i++) {
    $srcRow = 150 + $i
    $dstRow = 144 + $i
    $ws.Cells.Item($dstRow, 1).Value2 = $ws.Cells.Item($srcRow, 1).Value2
    $ws.Cells.Item($dstRow, 2).Value2 = $ws.Cells.Item($srcRow, 2).Value2
    $ws.Cells.Item($dstRow, 3).Value2 = $ws.Cells.Item($srcRow, 3).Value2
    $dVal = $ws.Cells.Item($srcRow, 4).Value2
    if ($dVal -ne $null) {
        $ws.Cells.Item($dstRow, 4).Value2 = $dVal
    }
}

# Step 3: remove the now-duplicated rows (old plane02 block, shifted to
# 150-153). This leaves the "cjklee" rows (untouched all along) sitting at
# 148-149, exactly where the target layout wants them.
$ws.Rows("150:153").Delete(-4162)

# Step 4: append the 14 new font-substitution rows (150-163).
$ws.Range("A150").Value2 = "otf-cjkaishob-h"
$ws.Range("B150").Value2 = "Identity-H"
$ws.Range("C150").Value2 = "DFKaiShoPro5-W5.otf"

$ws.Range("A151").Value2 = "otf-cjkaishob-v"
$ws.Range("B151").Value2 = "Identity-V"
$ws.Range("C151").Value2 = "DFKaiShoPro5-W5.otf"

$ws.Range("A152").Value2 = "otf-cjkaishoe-h"
$ws.Range("B152").Value2 = "Identity-H"
$ws.Range("C152").Value2 = "DFKaiShoPro5-W7.otf"

$ws.Range("A153").Value2 = "otf-cjkaishoe-v"
$ws.Range("B153").Value2 = "Identity-V"
$ws.Range("C153").Value2 = "DFKaiShoPro5-W7.otf"

$ws.Range("A154").Value2 = "otf-cjkaishor-h"
$ws.Range("B154").Value2 = "Identity-H"
$ws.Range("C154").Value2 = "DFKaiShoPro5-W5.otf"

$ws.Range("A155").Value2 = "otf-cjkaishor-v"
$ws.Range("B155").Value2 = "Identity-V"
$ws.Range("C155").Value2 = "DFKaiShoPro5-W5.otf"

$ws.Range("A156").Value2 = "otf-cjtuskub-h"
$ws.Range("B156").Value2 = "Identity-H"
$ws.Range("C156").Value2 = "FOT-TsukuAOldMinPr6N-B.otf"

$ws.Range("A157").Value2 = "otf-cjtuskub-v"
$ws.Range("B157").Value2 = "Identity-V"
$ws.Range("C157").Value2 = "FOT-TsukuAOldMinPr6N-B.otf"

$ws.Range("A158").Value2 = "otf-cjtuskul-h"
$ws.Range("B158").Value2 = "Identity-H"
$ws.Range("C158").Value2 = "FOT-TsukuAOldMinPr6N-L.otf"

$ws.Range("A159").Value2 = "otf-cjtuskul-v"
$ws.Range("B159").Value2 = "Identity-V"
$ws.Range("C159").Value2 = "FOT-TsukuAOldMinPr6N-L.otf"

$ws.Range("A160").Value2 = "otf-cjtuskur-h"
$ws.Range("B160").Value2 = "Identity-H"
$ws.Range("C160").Value2 = "FOT-TsukuAOldMinPr6N-R.otf"

$ws.Range("A161").Value2 = "otf-cjtuskur-v"
$ws.Range("B161").Value2 = "Identity-V"
$ws.Range("C161").Value2 = "FOT-TsukuAOldMinPr6N-R.otf"

$ws.Range("A162").Value2 = "otf-cjudmarur-h"
$ws.Range("B162").Value2 = "Identity-H"
$ws.Range("C162").Value2 = "A-OTF-UDShinMGoPr6N-Regular.otf"

$ws.Range("A163").Value2 = "otf-cjudmarur-v"
$ws.Range("B163").Value2 = "Identity-V"
$ws.Range("C163").Value2 = "A-OTF-UDShinMGoPr6N-Regular.otf"

# Step 5: update the view state (selection / scroll position) to match the
# edited area.
$excel.ActiveWindow.ScrollRow = 133
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D153").Select()
